$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the existing header style (from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add new data value in H2 (plain/default style, like the rest of row 2's numeric cells)
$ws.Range("H2").Value = 0
